# Upload all data and code for L&O Submission 5/16/24 attempt 2 at 13:42

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Data edits (columns I and J, rows 2-9) ---
$ws.Range("I2").Value = 3.5
$ws.Range("J2").Value = 0.6

$ws.Range("J3").Value = 0.1

$ws.Range("I4").Value = 3.5
$ws.Range("J4").Value = 0.6

$ws.Range("J5").Value = 0.1
$ws.Range("J6").Value = 0.1
$ws.Range("J7").Value = 0.1
$ws.Range("J8").Value = 0.1
$ws.Range("J9").Value = 0.1

# --- Selection / active cell change ---
$ws.Range("J2").Select()

# --- Window view tweak (position/size of the workbook window) ---
$excel.ActiveWindow.Top = 460
$excel.ActiveWindow.Height = 17540
$excel.ActiveWindow.Left = 0
$excel.ActiveWindow.Width = 28800
